$wb = $excel.ActiveWorkbook

# The handoff for 57c10ec3-...md failed: the "Ready for handoff" status is
# replaced everywhere by "Handoff transform failed", the (never produced)
# handoff file/hyperlink and handoff datetime are cleared, and the handoff
# reason flips from "Include" to "Ignored" on both locale sheets.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column
    $ws.Range("B2").Value = "Handoff transform failed"

    # Drop the "Latest Handoff File" hyperlink + its cell content (C2) -
    # no handoff file was produced.
    foreach ($hl in @($ws.Hyperlinks)) {
        if ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 3) {
            $hl.Delete()
        }
    }
    $ws.Range("C2").Clear()

    # No handoff happened, so the handoff datetime reverts to the "never"
    # sentinel value (same one already used for the handback columns).
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # The handoff reason is now "Ignored" instead of "Include".
    $ws.Range("H2").Value = "Ignored"
}

Write-Output "done"
